# Convert the raw "confession" submission entries into the anonymized
# "confession N" placeholders used for the posted images, and touch up the
# header row's height to match the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the two confession text entries.
$ws.Range("B2").Value = "confession 1"
$ws.Range("B3").Value = "confession 2"

# A2/A3 (the timestamp cells) pick up the exact same formatting as A1 --
# copy A1's format down onto them (xlPasteFormats = -4122).
$ws.Range("A1").Copy()
$ws.Range("A2:A3").PasteSpecial(-4122)

# Header row grows very slightly (18.75 -> 19.5), matching the data rows.
$ws.Rows.Item(1).RowHeight = 19.5
